# Add a second sheet ("naive") with NAIVE-model predictions alongside the
# existing sibling-regression sheet, so the two can be compared.
#
# 1. Rename the original "Sheet1" -> "sibregsimple"
# 2. Add a new worksheet "naive" after it, populate it with the NAIVE model
#    data (runyear / predicted return / p25 / p75), and turn it into a table
#    (Table13) styled like the original (TableStyleLight1) with the same
#    number formats (0 for runyear/p25/p75, 0.00 for predicted return).
# 3. Leave selection/active-sheet state matching the edited workbook:
#    sibregsimple shows A1:D15 selected (not the active tab any more),
#    naive is the active tab with F18 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# --- 1. rename the original sheet -----------------------------------------
$ws1.Name = "sibregsimple"

# --- 2. create + populate the new "naive" sheet ----------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "naive"

$headers = @("runyear", "predicted return", "p25", "p75")
for ($c = 1; $c -le 4; $c++) {
    $ws2.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$data = @(
    @(2012, 342174, 287001, 392441),
    @(2013, 315515, 293819, 342174),
    @(2014, 351087, 316834, 381393),
    @(2015, 339117, 298048, 380526),
    @(2016, 372239, 324549, 415229),
    @(2017, 297877, 234712, 369710),
    @(2018, 263094, 185352, 324013),
    @(2019, 167357, 141577, 187773),
    @(2020, 153755, 133983, 172805),
    @(2021, 148960, 134802, 164863),
    @(2022, 233795, 171847, 297664),
    @(2023, 348223, 267359, 427068),
    @(2024, 478498, 426099, 523874),
    @(2025, 593257, 534884, 654729)
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le 4; $c++) {
        $ws2.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

$lastRow = $r - 1

# number formats: runyear/p25/p75 -> "0", predicted return -> "0.00"
$ws2.Range("A2:A$lastRow").NumberFormat = "0"
$ws2.Range("B2:B$lastRow").NumberFormat = "0.00"
$ws2.Range("C2:D$lastRow").NumberFormat = "0"

# turn the range into a table matching the style of the original
$tableRange = $ws2.Range("A1:D$lastRow")
$lo2 = $ws2.ListObjects.Add(1, $tableRange, $null, 1)
$lo2.Name = "Table13"
$lo2.TableStyle = "TableStyleLight1"

# --- 3. selection / active sheet state -------------------------------------
$ws1.Range("A1:D15").Select() | Out-Null
$ws2.Range("F18").Select() | Out-Null
